$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: new resale-number record for 2024-01-07 00:09:47 (Sunday, week 01)
# Force text format on Date/Week columns so Excel doesn't auto-convert
# the literal strings ("2024-01-07" / "01") into a date serial or number.
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "2024-01-07"
$ws.Range("B29").Value = "00:09:47"
$ws.Range("C29").Value = "Sunday"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "01"

$ws.Range("E29").Value = 140532
$ws.Range("F29").Value = 143049
$ws.Range("G29").Value = 172337
$ws.Range("H29").Value = 147264
$ws.Range("I29").Value = -1
$ws.Range("J29").Value = 118396
$ws.Range("K29").Value = 224629
$ws.Range("L29").Value = 249366
$ws.Range("M29").Value = 185190
$ws.Range("N29").Value = 110410
$ws.Range("O29").Value = 40638
$ws.Range("P29").Value = 30808
$ws.Range("Q29").Value = 72518
$ws.Range("R29").Value = -1
$ws.Range("S29").Value = 42280
$ws.Range("T29").Value = -1
